$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET (sheet3) ---
$ws3 = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: name + week number
$ws3.Range("C1").Value = "Jesse Hare"
$ws3.Range("E1").Value = 11

# Row 3 - Stage/Task entries (set A3 and A4..A6 before B3 so shared-string
# insertion order matches: Jesse Hare, Project Design, Continue documentation,
# Code refactoring and optimisation, Testing, Implemetation of new requirements...)
$ws3.Range("A3").Value = "Project Design"

$ws3.Range("A4").Value = "Project Design"
$ws3.Range("B4").Value = "Continue documentation"
$ws3.Range("C4").Value = 4
$ws3.Range("D4").Value = 4
$ws3.Range("E4").Value = 0

$ws3.Range("A5").Value = "Project Design"
$ws3.Range("B5").Value = "Code refactoring and optimisation"
$ws3.Range("C5").Value = 3
$ws3.Range("D5").Value = 5
$ws3.Range("E5").Value = 0

$ws3.Range("A6").Value = "Project Design"
$ws3.Range("B6").Value = "Testing"
$ws3.Range("C6").Value = 2
$ws3.Range("D6").Value = 3
$ws3.Range("E6").Value = 0

$ws3.Range("B3").Value = "Implemetation of new requirements into program"
$ws3.Range("C3").Value = 8
$ws3.Range("D3").Value = 8
$ws3.Range("E3").Value = 0

# --- ACTIVITY LOG SUMMARY SHEET (sheet4) ---
$ws4 = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$ws4.Range("D1").Value = "Jesse Hare"

$ws4.Range("A4").Value = "Project Design"
$ws4.Range("B4").Value = 17
$ws4.Range("C4").Value = 3

# --- Active sheet / selection ---
$ws4.Activate()
$ws4.Range("D6").Select()
$ws3.Activate()
$ws3.Range("D10").Select()
